$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 8380.4
$ws.Range("J43").Value = 11967.333
$ws.Range("L43").Value = 11967.333
$ws.Range("N43").Value = -12105.333

$ws.Range("H74").Value = 33336000
$ws.Range("I74").Value = 33336000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 33336000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -33335064
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 33336000
$ws.Range("I77").Value = 33336000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 166680000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -166675320
$ws.Range("N77").ClearContents()

$ws.Range("H107").Value = 1147.5883
$ws.Range("I107").Value = 1175.25
$ws.Range("J107").Value = 1081.2
$ws.Range("K107").Value = 1175.25
$ws.Range("L107").Value = 1081.2
$ws.Range("M107").Value = 744.75
$ws.Range("N107").Value = -4921.2

$ws.Range("H116").Value = 1120622.8
$ws.Range("I116").Value = 2502649.8
$ws.Range("J116").Value = 15001.2
$ws.Range("K116").Value = 2502649.8
$ws.Range("L116").Value = 15001.2
$ws.Range("M116").Value = -2499207.8
$ws.Range("N116").Value = -21885.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7637.96
$ws.Range("I32").Value = 4737.7646
$ws.Range("J32").Value = 13800.875
$ws.Range("K32").Value = 4737.7646
$ws.Range("L32").Value = 13800.875
$ws.Range("M32").Value = -4450.7646
$ws.Range("N32").Value = -14374.875

$ws.Range("H51").Value = 51588.75
$ws.Range("J51").Value = 51588.75
$ws.Range("L51").Value = 51588.75
$ws.Range("N51").Value = -53100.75

$ws.Range("H97").Value = 1120.0555
$ws.Range("I97").Value = 957.4
$ws.Range("K97").Value = 957.4
$ws.Range("M97").Value = -461.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 4439.8
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 5424.75
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 5424.75
$ws.Range("M10").Value = -360
$ws.Range("N10").Value = -5704.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2144.2222
$ws.Range("I5").Value = 319.6
$ws.Range("J5").Value = 4425
$ws.Range("K5").Value = 319.6
$ws.Range("L5").Value = 4425
$ws.Range("M5").Value = -207.6
$ws.Range("N5").Value = -4649

$ws.Range("H31").Value = 3990.9443
$ws.Range("I31").Value = 1602.8462
$ws.Range("K31").Value = 1602.8462
$ws.Range("M31").Value = -1307.8462

$ws.Range("H34").Value = 3990.9443
$ws.Range("I34").Value = 1602.8462
$ws.Range("K34").Value = 1602.8462
$ws.Range("M34").Value = -1400.8462

$ws.Range("H122").Value = 5816.8335
$ws.Range("I122").Value = 4304
$ws.Range("J122").Value = 7329.6665
$ws.Range("K122").Value = 12912
$ws.Range("L122").Value = 21988.9995
$ws.Range("M122").Value = -10462
$ws.Range("N122").Value = -26888.9995

$ws.Range("H137").Value = 49780
$ws.Range("J137").Value = 49780
$ws.Range("L137").Value = 49780
$ws.Range("N137").Value = -59980

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 12006.158
$ws.Range("I34").Value = 22337.8
$ws.Range("J34").Value = 8316.286
$ws.Range("K34").Value = 67013.4
$ws.Range("L34").Value = 24948.858
$ws.Range("M34").Value = -66929.4
$ws.Range("N34").Value = -25116.858

$ws.Range("H39").Value = 13617.333
$ws.Range("J39").Value = 13477.177
$ws.Range("L39").Value = 40431.531
$ws.Range("N39").Value = -41019.531

$ws.Range("H75").Value = 1212
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 1212
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H87").Value = 6352
$ws.Range("I87").Value = 6352
$ws.Range("K87").Value = 19056
$ws.Range("M87").Value = -17808

$ws.Range("H90").Value = 6352
$ws.Range("I90").Value = 6352
$ws.Range("K90").Value = 57168
$ws.Range("M90").Value = -50928

$ws.Range("H107").Value = 1305.4117
$ws.Range("I107").Value = 440
$ws.Range("J107").Value = 1666
$ws.Range("K107").Value = 1320
$ws.Range("L107").Value = 4998
$ws.Range("M107").Value = 600
$ws.Range("N107").Value = -8838

$ws.Range("H113").Value = 615.59375
$ws.Range("I113").Value = 653.5
$ws.Range("J113").Value = 586.1111
$ws.Range("K113").Value = 1960.5
$ws.Range("L113").Value = 1758.3333
$ws.Range("M113").Value = 209.5
$ws.Range("N113").Value = -6098.3333

$ws.Range("H114").Value = 41667884
$ws.Range("I114").Value = 90909280
$ws.Range("J114").Value = 2090.3845
$ws.Range("K114").Value = 272727840
$ws.Range("L114").Value = 6271.1535
$ws.Range("M114").Value = -272724586
$ws.Range("N114").Value = -12779.1535

$ws.Range("H117").Value = 1099.125
$ws.Range("I117").Value = 543
$ws.Range("J117").Value = 1432.8
$ws.Range("K117").Value = 1629
$ws.Range("L117").Value = 4298.4
$ws.Range("M117").Value = 1813
$ws.Range("N117").Value = -11182.4

$ws.Range("H121").Value = 2107.776
$ws.Range("I121").Value = 440
$ws.Range("J121").Value = 2198.7454
$ws.Range("K121").Value = 1320
$ws.Range("L121").Value = 6596.236199999999
$ws.Range("M121").Value = -10
$ws.Range("N121").Value = -9216.2362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10418838
$ws.Range("I22").Value = 27779978
$ws.Range("J22").Value = 2153.6667
$ws.Range("K22").Value = 27779978
$ws.Range("L22").Value = 2153.6667
$ws.Range("M22").Value = -27779683
$ws.Range("N22").Value = -2743.6667

$ws.Range("H27").Value = 10418838
$ws.Range("I27").Value = 27779978
$ws.Range("J27").Value = 2153.6667
$ws.Range("K27").Value = 27779978
$ws.Range("L27").Value = 2153.6667
$ws.Range("M27").Value = -27779871
$ws.Range("N27").Value = -2367.6667

$ws.Range("H46").Value = 2416.6667
$ws.Range("I46").Value = 2837.5
$ws.Range("J46").Value = 2263.6365
$ws.Range("K46").Value = 2837.5
$ws.Range("L46").Value = 2263.6365
$ws.Range("M46").Value = -2649.5
$ws.Range("N46").Value = -2639.6365

$ws.Range("H61").Value = 1606.7273
$ws.Range("I61").Value = 1680.5714
$ws.Range("J61").Value = 1477.5
$ws.Range("K61").Value = 1680.5714
$ws.Range("L61").Value = 1477.5
$ws.Range("M61").Value = -1478.5714
$ws.Range("N61").Value = -1881.5

$ws.Range("H113").Value = 1606.7273
$ws.Range("I113").Value = 1680.5714
$ws.Range("J113").Value = 1477.5
$ws.Range("K113").Value = 1680.5714
$ws.Range("L113").Value = 1477.5
$ws.Range("M113").Value = 489.4286
$ws.Range("N113").Value = -5817.5

$ws.Range("H132").Value = 6072.815
$ws.Range("I132").Value = 1520.8
$ws.Range("J132").Value = 8750.471
$ws.Range("K132").Value = 4562.4
$ws.Range("L132").Value = 26251.413
$ws.Range("M132").Value = -2032.4
$ws.Range("N132").Value = -31311.413

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 25374.5
$ws.Range("I113").Value = 100000
$ws.Range("K113").Value = 300000
$ws.Range("M113").Value = -297830

$ws.Range("H122").Value = 7259
$ws.Range("I122").Value = 6945
$ws.Range("K122").Value = 20835
$ws.Range("M122").Value = -18385

$ws.Range("H132").Value = 12351772
$ws.Range("I132").Value = 9982.818
$ws.Range("J132").Value = 20836752
$ws.Range("K132").Value = 29948.454
$ws.Range("L132").Value = 62510256
$ws.Range("M132").Value = -27418.454
$ws.Range("N132").Value = -62515316

$ws.Range("H136").Value = 4119.4165
$ws.Range("I136").Value = 1162.3334
$ws.Range("J136").Value = 7076.5
$ws.Range("K136").Value = 3487.0002
$ws.Range("L136").Value = 21229.5
$ws.Range("M136").Value = -937.0002
$ws.Range("N136").Value = -26329.5
